$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Update the "time_taken" (column F) timestamps on the "data" sheet for
#    rows 2..113 to the refreshed query timestamps.
# ---------------------------------------------------------------------------
$timestamps = @{
    2 = "2021-10-05 14:19:42.021913"
    3 = "2021-10-05 14:19:42.021921"
    4 = "2021-10-05 14:19:42.021924"
    5 = "2021-10-05 14:19:42.021927"
    6 = "2021-10-05 14:19:42.021930"
    7 = "2021-10-05 14:19:42.021932"
    8 = "2021-10-05 14:19:42.021935"
    9 = "2021-10-05 14:19:42.021937"
    10 = "2021-10-05 14:19:42.021940"
    11 = "2021-10-05 14:19:42.021943"
    12 = "2021-10-05 14:19:42.021945"
    13 = "2021-10-05 14:19:42.021948"
    14 = "2021-10-05 14:19:42.021950"
    15 = "2021-10-05 14:19:42.021953"
    16 = "2021-10-05 14:19:42.021955"
    17 = "2021-10-05 14:19:42.021958"
    18 = "2021-10-05 14:19:42.021961"
    19 = "2021-10-05 14:19:42.021963"
    20 = "2021-10-05 14:19:42.021966"
    21 = "2021-10-05 14:19:42.021968"
    22 = "2021-10-05 14:19:42.021971"
    23 = "2021-10-05 14:19:42.021973"
    24 = "2021-10-05 14:19:42.021976"
    25 = "2021-10-05 14:19:42.021978"
    26 = "2021-10-05 14:19:42.021981"
    27 = "2021-10-05 14:19:42.021984"
    28 = "2021-10-05 14:19:42.021986"
    29 = "2021-10-05 14:19:42.021989"
    30 = "2021-10-05 14:19:42.021991"
    31 = "2021-10-05 14:19:42.021994"
    32 = "2021-10-05 14:19:42.021996"
    33 = "2021-10-05 14:19:42.021999"
    34 = "2021-10-05 14:19:42.022001"
    35 = "2021-10-05 14:19:42.022004"
    36 = "2021-10-05 14:19:42.022007"
    37 = "2021-10-05 14:19:42.022009"
    38 = "2021-10-05 14:19:42.022012"
    39 = "2021-10-05 14:19:42.022014"
    40 = "2021-10-05 14:19:42.022017"
    41 = "2021-10-05 14:19:42.022019"
    42 = "2021-10-05 14:19:42.022022"
    43 = "2021-10-05 14:19:42.022025"
    44 = "2021-10-05 14:19:42.022027"
    45 = "2021-10-05 14:19:42.022030"
    46 = "2021-10-05 14:19:42.022032"
    47 = "2021-10-05 14:19:42.022035"
    48 = "2021-10-05 14:19:42.022037"
    49 = "2021-10-05 14:19:42.022040"
    50 = "2021-10-05 14:19:42.022042"
    51 = "2021-10-05 14:19:42.022045"
    52 = "2021-10-05 14:19:42.022047"
    53 = "2021-10-05 14:19:42.022050"
    54 = "2021-10-05 14:19:42.022053"
    55 = "2021-10-05 14:19:42.022055"
    56 = "2021-10-05 14:19:42.022058"
    57 = "2021-10-05 14:19:42.022060"
    58 = "2021-10-05 14:19:42.022063"
    59 = "2021-10-05 14:19:42.022065"
    60 = "2021-10-05 14:19:42.022068"
    61 = "2021-10-05 14:19:42.022070"
    62 = "2021-10-05 14:19:42.022073"
    63 = "2021-10-05 14:19:42.022075"
    64 = "2021-10-05 14:19:42.022078"
    65 = "2021-10-05 14:19:42.022080"
    66 = "2021-10-05 14:19:42.022084"
    67 = "2021-10-05 14:19:42.022087"
    68 = "2021-10-05 14:19:42.022089"
    69 = "2021-10-05 14:19:42.022092"
    70 = "2021-10-05 14:19:42.022094"
    71 = "2021-10-05 14:19:42.022097"
    72 = "2021-10-05 14:19:42.022099"
    73 = "2021-10-05 14:19:42.022102"
    74 = "2021-10-05 14:19:42.022104"
    75 = "2021-10-05 14:19:42.022107"
    76 = "2021-10-05 14:19:42.022109"
    77 = "2021-10-05 14:19:42.022112"
    78 = "2021-10-05 14:19:42.022116"
    79 = "2021-10-05 14:19:42.022119"
    80 = "2021-10-05 14:19:42.022122"
    81 = "2021-10-05 14:19:42.022124"
    82 = "2021-10-05 14:19:42.022127"
    83 = "2021-10-05 14:19:42.022129"
    84 = "2021-10-05 14:19:42.022132"
    85 = "2021-10-05 14:19:42.022135"
    86 = "2021-10-05 14:19:42.022137"
    87 = "2021-10-05 14:19:42.022140"
    88 = "2021-10-05 14:19:42.022142"
    89 = "2021-10-05 14:19:42.022145"
    90 = "2021-10-05 14:19:42.022147"
    91 = "2021-10-05 14:19:42.022150"
    92 = "2021-10-05 14:19:42.022152"
    93 = "2021-10-05 14:19:42.022155"
    94 = "2021-10-05 14:19:42.022158"
    95 = "2021-10-05 14:19:42.022161"
    96 = "2021-10-05 14:19:42.022164"
    97 = "2021-10-05 14:19:42.022166"
    98 = "2021-10-05 14:19:42.022169"
    99 = "2021-10-05 14:19:42.022171"
    100 = "2021-10-05 14:19:42.022174"
    101 = "2021-10-05 14:19:42.022176"
    102 = "2021-10-05 14:19:42.022179"
    103 = "2021-10-05 14:19:42.022182"
    104 = "2021-10-05 14:19:42.022184"
    105 = "2021-10-05 14:19:42.022187"
    106 = "2021-10-05 14:19:42.022189"
    107 = "2021-10-05 14:19:42.022192"
    108 = "2021-10-05 14:19:42.022194"
    109 = "2021-10-05 14:19:42.022197"
    110 = "2021-10-05 14:19:42.022202"
    111 = "2021-10-05 14:19:42.022205"
    112 = "2021-10-05 14:19:42.022207"
    113 = "2021-10-05 14:19:42.022210"
}
foreach ($row in $timestamps.Keys) {
    $dataSheet.Cells.Item([int]$row, 6).Value = $timestamps[$row]
}

# ---------------------------------------------------------------------------
# 2. Add a new "metadata" worksheet right after "data".
# ---------------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (B1:G1) -- copy the bold/bordered header style used on "data"
# (columns B1:F1), then stamp the extra G1 header cell with the same style
# copied from a single source cell so no new style is introduced.
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("B1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row 2 -- A2 mirrors the styled index column used throughout "data".
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$metaSheet.Range("A2").Value = 0

$metaSheet.Range("B2").Value = "Congenital disorders of glycosylation"
$metaSheet.Range("C2").Value = 25

# data_version must stay textual ("2.76"), not become a numeric 2.76.
$metaSheet.Range("D2").Value = "'2.76"
$dataSheet.Range("B2").Copy()
$metaSheet.Range("D2").PasteSpecial(-4122)

$metaSheet.Range("E2").Value = "2021-09-02T17:02:45.175551Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:19:42.018206"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/25/?format=json"

$dataSheet.Select()

Write-Host "metadata sheet added and timestamps refreshed"
